$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up physical column layout so the surviving custom column widths
# (old H/I => new G/H) shift left by one, and a fresh column is opened at I
# for the new "AREA" field (mirrors what the author's edit did: one column
# removed from the B..H block, one new column opened before "CREADO POR").
$ws.Columns.Item(3).Delete()
$ws.Columns.Item(9).Insert()

# New column I needs its own explicit width (it did not exist before).
$ws.Columns.Item(9).ColumnWidth = 15.5703125

# --- Give row 11 the same row-level formatting as row 10 before filling it in.
$ws.Range("A10:K10").Copy()
$ws.Range("A11:K11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row ---
$ws.Range("A1").Value = "LUGAR"
$ws.Range("B1").Value = "CAMPUS"
$ws.Range("C1").Value = "EDIFICIO"
$ws.Range("D1").Value = "ESPACIO"
$ws.Range("E1").Value = "PLANTA"
$ws.Range("F1").Value = "ESTANCIA"
$ws.Range("G1").Value = "ELEMENTOS DEL ESPACIO"
$ws.Range("H1").Value = "DESCRIPCIÓN"
$ws.Range("I1").Value = "AREA"
$ws.Range("J1").Value = "CREADO POR"
$ws.Range("K1").Value = "SUPERVISOR"

# --- Data rows 2-10 ---
$rooms = @("HALL","DESPACHO 1","DESPACHO DERECHA","ESCALERA","ASCENSORES","OFICINA 1","OFICINA 1","OFICINA 2","ASCENSORES")
$plantas = @(0,0,0,0,0,1,1,1,1)
for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $ws.Range("A$r").Value = "MADRID"
    $ws.Range("B$r").Value = "LAS TABLAS"
    $ws.Range("C$r").Value = "EDF. 4"
    $ws.Range("D$r").Value = ""
    $ws.Range("E$r").Value = $plantas[$i]
    $ws.Range("F$r").Value = $rooms[$i]
    $ws.Range("G$r").Value = ""
    $ws.Range("H$r").Value = ""
    $ws.Range("I$r").Value = ""
    $ws.Range("J$r").Value = "CARLOS"
    $ws.Range("K$r").Value = "CARLOS"
}
# HALL row carries a description.
$ws.Range("H2").Value = "El hall principal, suelo de mármol"

# --- Row 11: new "PARQUE" space entry ---
$ws.Range("A11").Value = "MADRID"
$ws.Range("B11").Value = "LAS TABLAS"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = "PARQUE"
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "Es un parque"
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = "CARLOS"
$ws.Range("K11").Value = "CARLOS"

# --- Stray formatted-but-empty cells left over from typing + deleting text
# (underlined default font, no fill/border) ---
$ws.Range("J14").Font.Underline = $true
$ws.Range("F16").Font.Underline = $true

# --- Conditional formatting now covers the new row too ---
$ws.Range("A2:K9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("A2:K11"))

# --- The old "Floors and Rooms / Rooms Only" list validation no longer applies ---
$ws.Range("E2:E10").Validation.Delete()

# --- Selection left where the user was last working ---
$ws.Range("F15").Select()

Write-Output "done"
